$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 metrics (Appenzeller-Herzog 2019) per corrected relevance markers
$ws.Range("H3").Value = 0.7876783849634529
$ws.Range("I3").Value = 0.06610196968468833
$ws.Range("K3").Value = 201.6923076923077

$ws.Range("Q3").Value = 46
$ws.Range("R3").Value = 54
$ws.Range("S3").Value = 89
$ws.Range("T3").Value = 144
$ws.Range("U3").Value = 258

$ws.Range("V3").Value = 2801
$ws.Range("W3").Value = 2793
$ws.Range("X3").Value = 2758
$ws.Range("Y3").Value = 2703
$ws.Range("Z3").Value = 2589

$ws.Range("AF3").Value = 0.983843
$ws.Range("AG3").Value = 0.981033
$ws.Range("AH3").Value = 0.968739
$ws.Range("AI3").Value = 0.94942
$ws.Range("AJ3").Value = 0.909378
